$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordered data: label, count
$labels = @(
  "Wrong_Entity_Event_as_NonEvent",
  "Wrong_Entity_NonEvent_as_Event",
  "Correct",
  "Wrong_Tag_B_as_I",
  "Wrong_Tag_E_as_I",
  "Wrong_Tag_B_as_E"
)
$values = @(77, 69, 48, 4, 2, 1)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Ensure the newly added row (row 7) column A cell carries the same style
# (bold/centered/bordered) as the other label cells in column A.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(7, 1).Value = $labels[5]

